$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new status row (row 22) below the existing data, reusing the
# same date style (numFmtId 14, "m/d/yyyy") already used by column A.
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = 43196

$ws.Range("B22").Value = "Android Layout (20%)"
$ws.Range("C22").Value = "C# small layout adaption"
$ws.Range("D22").Value = "WebApp layout (25%)"

$ws.Range("D26").Select()
